# Regenerate the "K" (strikeouts) column (G) of the game log using the
# correct source values instead of the old "Strike#" values, and refresh a
# couple of dependent cells (H2/I2) that were recomputed alongside it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G ("K") values, row by row (row 2 .. row 39)
$newK = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 0
    10 = 0
    11 = 0
    12 = 1
    13 = 0
    14 = 1
    15 = 0
    16 = 3
    17 = 1
    18 = 2
    19 = 2
    20 = 2
    21 = 1
    22 = 0
    23 = 2
    24 = 1
    25 = 1
    26 = 1
    27 = 5
    28 = 3
    29 = 3
    30 = 1
    31 = 1
    32 = 0
    33 = 2
    34 = 2
    35 = 2
    36 = 1
    37 = 0
    38 = 3
    39 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}

# Dependent recalculated values on row 2 (IP and I0) that shifted alongside
# the corrected K value.
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 6
